$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2009 data row (row 2) is removed; subsequent rows (2010, 2011) shift
# up by one, so the old row 3 becomes row 2 and old row 4 becomes row 3.
$ws.Rows("2:2").Delete()
